$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.451.57'
$ws.Range('E2').Value = '  +1.50%  '

$ws.Range('D3').Value = '3.006.50'
$ws.Range('E3').Value = '  -0.03%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '508.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.10%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('E8').Value = '  +1.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.50'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.42%  '

$ws.Range('E10').Value = '  +1.19%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.364'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.37%  '

$ws.Range('D12').Value = '3.521.62'
$ws.Range('E12').Value = '  -0.02%  '

$ws.Range('E13').Value = '  +0.84%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.38'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.60%  '

$ws.Range('E15').Value = '  +5.83%  '

$ws.Range('D16').Value = '57.479.13'

$ws.Range('E17').Value = '  +6.77%  '

$ws.Range('D18').Value = '3.008.27'
$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.02%  '

$ws.Range('E25').Value = '  -0.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.15%  '

$ws.Range('E27').Value = '  +0.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.24%  '

$ws.Range('E29').Value = '  +4.65%  '

$ws.Range('E30').Value = '  +1.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.18%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.56'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.73%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.46%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '153.73'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.47%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.85'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.85%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.44%  '

$ws.Range('E37').Value = '  +0.92%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.37'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.82%  '

$ws.Range('D39').Value = '3.041.16'
$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('E40').Value = '  +1.74%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.84'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.10%  '

$ws.Range('E43').Value = '  -0.30%  '

$ws.Range('D44').Value = '2.266.18'
$ws.Range('E44').Value = '  -0.59%  '

$ws.Range('E45').Value = '  -0.29%  '

$ws.Range('E46').Value = '  -1.72%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.01'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.16%  '

$ws.Range('E48').Value = '  +1.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.28%  '

$ws.Range('E50').Value = '  -7.16%  '

$ws.Range('E51').Value = '  +2.43%  '
